# Fruta / hortaliza, semanal
# A new weekly price observation is inserted as row 129 (date 2022-12-29),
# pushing the previously existing rows 129-172 down to rows 130-173.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above what is currently row 129. This shifts the
# existing rows 129-172 down to 130-173 and extends the used range /
# dimension to A1:R173 automatically.
$ws.Rows.Item(129).Insert()

# Populate the newly inserted row 129 with the new weekly record.
$ws.Cells.Item(129, 1).Value = 11
$ws.Cells.Item(129, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(129, 3).Value = "Bíobío"
$ws.Cells.Item(129, 4).Value = 44924
$ws.Cells.Item(129, 5).Value = 8
$ws.Cells.Item(129, 6).Value = 100112043
$ws.Cells.Item(129, 7).Value = "Pepino ensalada"
$ws.Cells.Item(129, 8).Value = "Sin especificar"
$ws.Cells.Item(129, 9).Value = "Primera"
$ws.Cells.Item(129, 10).Value = 100
$ws.Cells.Item(129, 11).Value = 16000
$ws.Cells.Item(129, 12).Value = 17000
$ws.Cells.Item(129, 13).Value = 16500
$ws.Cells.Item(129, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(129, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(129, 16).Value = 275
$ws.Cells.Item(129, 17).Value = 60
$ws.Cells.Item(129, 18).Value = "Hortaliza"
